# Auto-generated Excel COM-interop script
# Applies updated market/profit figures (columns H-N) for the rows
# touched by the scheduled market-data refresh, across all 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 33000
$ws.Range("J81").Value = 33000
$ws.Range("L81").Value = 33000
$ws.Range("N81").Value = -34996

$ws.Range("H84").Value = 33000
$ws.Range("J84").Value = 33000
$ws.Range("L84").Value = 99000
$ws.Range("N84").Value = -108984

$ws.Range("H96").Value = 1575.5
$ws.Range("I96").Value = 524.8461
$ws.Range("J96").Value = 2817.182
$ws.Range("K96").Value = 1574.5383
$ws.Range("L96").Value = 8451.545999999998
$ws.Range("M96").Value = -201.5382999999999
$ws.Range("N96").Value = -11197.546

$ws.Range("H98").Value = 1010.6286
$ws.Range("I98").Value = 901.5161000000001
$ws.Range("J98").Value = 1856.25
$ws.Range("K98").Value = 901.5161000000001
$ws.Range("L98").Value = 1856.25
$ws.Range("M98").Value = 596.4838999999999
$ws.Range("N98").Value = -4852.25

$ws.Range("H105").Value = 52499
$ws.Range("J105").Value = 52499
$ws.Range("L105").Value = 52499
$ws.Range("N105").Value = -59487

$ws.Range("H112").Value = 2131.3076
$ws.Range("J112").Value = 2191.2
$ws.Range("L112").Value = 6573.599999999999
$ws.Range("N112").Value = -8789.599999999999

$ws.Range("H113").Value = 3110.5417
$ws.Range("I113").Value = 2939.85
$ws.Range("K113").Value = 2939.85
$ws.Range("M113").Value = 314.1500000000001

$ws.Range("H122").Value = 1010.6286
$ws.Range("I122").Value = 901.5161000000001
$ws.Range("J122").Value = 1856.25
$ws.Range("K122").Value = 2704.5483
$ws.Range("L122").Value = 5568.75
$ws.Range("M122").Value = -254.5483000000004
$ws.Range("N122").Value = -10468.75

$ws.Range("H141").Value = 13659.046
$ws.Range("I141").Value = 5166.3335
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 15499.0005
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -10319.0005
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1063077
$ws.Range("I32").Value = 1140009.1
$ws.Range("K32").Value = 1140009.1
$ws.Range("M32").Value = -1139722.1

$ws.Range("H61").Value = 3126904.2
$ws.Range("I61").Value = 2004.9286
$ws.Range("K61").Value = 2004.9286
$ws.Range("M61").Value = -1792.9286

$ws.Range("H122").Value = 1618
$ws.Range("I122").Value = 1632.75
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4898.25
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2448.25
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 4410.931
$ws.Range("J132").Value = 7441
$ws.Range("L132").Value = 22323
$ws.Range("N132").Value = -27383

$ws.Range("H136").Value = 3126904.2
$ws.Range("I136").Value = 2004.9286
$ws.Range("K136").Value = 6014.7858
$ws.Range("M136").Value = -3464.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5556906.5
$ws.Range("I134").Value = 1447.5
$ws.Range("K134").Value = 4342.5
$ws.Range("M134").Value = -1807.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4687.375
$ws.Range("I16").Value = 4350.25
$ws.Range("J16").Value = 5024.5
$ws.Range("K16").Value = 4350.25
$ws.Range("L16").Value = 5024.5
$ws.Range("M16").Value = -4063.25
$ws.Range("N16").Value = -5598.5

$ws.Range("H105").Value = 35511.223
$ws.Range("I105").Value = 35511.223
$ws.Range("K105").Value = 35511.223
$ws.Range("M105").Value = -33764.223

$ws.Range("H107").Value = 391.78262
$ws.Range("I107").Value = 291.86667
$ws.Range("K107").Value = 291.86667
$ws.Range("M107").Value = 1628.13333

$ws.Range("H113").Value = 4687.375
$ws.Range("I113").Value = 4350.25
$ws.Range("J113").Value = 5024.5
$ws.Range("K113").Value = 4350.25
$ws.Range("L113").Value = 5024.5
$ws.Range("M113").Value = -2180.25
$ws.Range("N113").Value = -9364.5

$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

$ws.Range("H132").Value = 1960.7709
$ws.Range("I132").Value = 1822.4878
$ws.Range("K132").Value = 5467.463400000001
$ws.Range("M132").Value = -2937.463400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 777.6667
$ws.Range("J113").Value = 802.9091
$ws.Range("L113").Value = 2408.7273
$ws.Range("N113").Value = -6748.7273

$ws.Range("H136").Value = 11250.417
$ws.Range("I136").Value = 6001.25
$ws.Range("K136").Value = 18003.75
$ws.Range("M136").Value = -12903.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 19000
$ws.Range("J53").Value = 19000
$ws.Range("L53").Value = 19000
$ws.Range("N53").Value = -20262

$ws.Range("H102").Value = 1178.8334
$ws.Range("I102").Value = 1014.6
$ws.Range("K102").Value = 1014.6
$ws.Range("M102").Value = 607.4

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws.Range("H132").Value = 15897.792
$ws.Range("I132").Value = 9086.556
$ws.Range("K132").Value = 27259.668
$ws.Range("M132").Value = -24729.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3600
$ws.Range("I7").Value = 3600
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3600
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3488
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 3047.1
$ws.Range("J22").Value = 3308.7273
$ws.Range("L22").Value = 3308.7273
$ws.Range("N22").Value = -3898.7273

$ws.Range("H27").Value = 3047.1
$ws.Range("J27").Value = 3308.7273
$ws.Range("L27").Value = 3308.7273
$ws.Range("N27").Value = -3522.7273

$ws.Range("H31").Value = 2421.4
$ws.Range("J31").Value = 3730.6667
$ws.Range("L31").Value = 3730.6667
$ws.Range("N31").Value = -4226.6667

$ws.Range("H46").Value = 4497.2
$ws.Range("I46").Value = 2139
$ws.Range("J46").Value = 9999.666999999999
$ws.Range("K46").Value = 2139
$ws.Range("L46").Value = 9999.666999999999
$ws.Range("M46").Value = -1951
$ws.Range("N46").Value = -10375.667

$ws.Range("H122").Value = 3380.45
$ws.Range("J122").Value = 4307.5
$ws.Range("L122").Value = 12922.5
$ws.Range("N122").Value = -17822.5

$ws.Range("H126").Value = 3600
$ws.Range("I126").Value = 3600
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10800
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8330
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 13160996
$ws.Range("I136").Value = 7815308
$ws.Range("J136").Value = 41671332
$ws.Range("K136").Value = 23445924
$ws.Range("L136").Value = 125013996
$ws.Range("M136").Value = -23443374
$ws.Range("N136").Value = -125019096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 17999
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 17999
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 17999
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -18451

$ws.Range("H100").Value = 455.3
$ws.Range("I100").Value = 455.3
$ws.Range("K100").Value = 910.6
$ws.Range("M100").Value = -369.6

$ws.Range("H126").Value = 6752
$ws.Range("I126").Value = 10004
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 30012
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -27542
$ws.Range("N126").Value = -15440

$ws.Range("H136").Value = 5351224
$ws.Range("I136").Value = 2900901.5
$ws.Range("J136").Value = 12702193
$ws.Range("K136").Value = 8702704.5
$ws.Range("L136").Value = 38106579
$ws.Range("M136").Value = -8700154.5
$ws.Range("N136").Value = -38111679

Write-Host "Applied scheduled market-data refresh to all sheets."
